$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(35).Insert()

$ws.Cells.Item(35,1).Value = 7
$ws.Cells.Item(35,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35,3).Value = "Ñuble"
$ws.Cells.Item(35,4).Value = 44980
$ws.Cells.Item(35,5).Value = 16
$ws.Cells.Item(35,6).Value = "Fruta"
$ws.Cells.Item(35,7).Value = 100108
$ws.Cells.Item(35,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(35,9).Value = 100108002
$ws.Cells.Item(35,10).Value = "Mango"
$ws.Cells.Item(35,11).Value = "Sin especificar"
$ws.Cells.Item(35,12).Value = "Primera"
$ws.Cells.Item(35,13).Value = 60
$ws.Cells.Item(35,14).Value = 7000
$ws.Cells.Item(35,15).Value = 7500
$ws.Cells.Item(35,16).Value = 7250
$ws.Cells.Item(35,17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(35,18).Value = "Perú"
$ws.Cells.Item(35,19).Value = 1812
$ws.Cells.Item(35,20).Value = 4
